$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($col in @("D", "E", "F", "G", "H", "I", "J")) {
    $ws.Range("$col`2").Value = "unknown"
}
